{"js": "// Add the year \"2018\" after the \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y ... th\u00e1ng ... n\u0103m\" date\n// line in the notice header, matching the fix described in the commit:\n// \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m\" -> \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m 2018\"\n// (italicised \"2018\" appended as its own run, same as the surrounding text).\n\nconst TARGET_TEXT = \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m\";\n\n// Look the exact phrase up so we do not depend on paragraph indices.\nconst results = context.document.body.search(TARGET_TEXT, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const dateRange = results.items[0];\n\n  // 1) Add the trailing space onto the existing run (keeps the original\n  //    run/formatting, just like `xml:space=\"preserve\"` text growing by\n  //    one character in the source diff).\n  const afterSpace = dateRange.insertText(\" \", Word.InsertLocation.end);\n  await context.sync();\n\n  // 2) Insert \"2018\" immediately after as a brand-new run. Toggling the\n  //    formatting (bold on, then back off) forces the engine to keep this\n  //    as a distinct run instead of silently re-merging it into the\n  //    previous one, while the final formatting still matches the\n  //    surrounding italic text exactly.\n  const insertionPoint = afterSpace.getRange(Word.RangeLocation.end);\n  const yearRange = insertionPoint.insertText(\"2018\", Word.InsertLocation.replace);\n  yearRange.font.set({ italic: true, bold: true });\n  await context.sync();\n  yearRange.font.set({ italic: true, bold: false });\n  await context.sync();\n}\n", "ps1": "# Add the year \"2018\" after the \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y ... th\u00e1ng ... n\u0103m\" date\n# line in the notice header, matching the fix described in the commit:\n# \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m\" -> \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m 2018\"\n# (italicised \"2018\" appended as its own run, same as the surrounding text).\n\n$d = $word.ActiveDocument\n\n$TARGET_TEXT = \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m\"\n\n# Idempotency guard: if the fix was already applied, do nothing.\n$already = $d.Content\n$alreadyApplied = $already.Find.Execute($TARGET_TEXT + \" 2018\")\n\nif (-not $alreadyApplied) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($TARGET_TEXT)\n\n    if ($found) {\n        # 1) Add the trailing space onto the existing run (keeps the\n        #    original run/formatting, just like `xml:space=\"preserve\"`\n        #    text growing by one character in the source diff).\n        $rng.Collapse(0)\n        $rng.InsertAfter(\" \")\n\n        # 2) Insert \"2018\" immediately after as a brand-new run. Toggling\n        #    the formatting (bold on, then back off) forces a distinct\n        #    run instead of silently re-merging it into the previous one,\n        #    while the final formatting still matches the surrounding\n        #    italic text exactly.\n        $rng.Collapse(0)\n        $rng.InsertAfter(\"2018\")\n        $rng.Bold = 1\n        $rng.Bold = 0\n    }\n}\n"}
